# Inserts a new weekly price record at row 16 (pushing the existing rows
# 16-23 down to 17-24) on the Camote/La Palmera de La Serena sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 16:23 down to 17:24, duplicating row 16's formatting (date style)
# onto the new row 16.
$ws.Rows("16:16").Insert()

# Populate the newly inserted row 16 with the latest week's data.
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = "Terminal La Palmera de La Serena"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 45229
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 100114002
$ws.Range("G16").Value = "Camote"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 460
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 17000
$ws.Range("M16").Value = 16500
$ws.Range("N16").Value = "`$/malla 18 kilos"
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 917
$ws.Range("Q16").Value = 18
$ws.Range("R16").Value = "Hortaliza"
